$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name in B3 (remove "Tunis Re Société anonyme" suffix)
$ws.Range("B3").Value = "Société Tunisienne de Réassurance (BVMT:TRE)"

# --- Row 2 ---
$ws.Range("D2").Value = 0.117
$ws.Range("E2").Value = 0.0236
$ws.Range("G2").Value = 0.06515151515151515
$ws.Range("H2").Value = 0.06515151515151515
$ws.Range("I2").Value = 0.1608585858585859
$ws.Range("J2").Value = 0.1254877203406615
$ws.Range("K2").Value = 5.57
$ws.Range("L2").Value = 0.1406565656565656
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 2.55
$ws.Range("V2").Value = 0.04419410745233968
$ws.Range("W2").Value = 0.08203240058910162
$ws.Range("X2").Value = 0.1351317619518854
$ws.Range("Y2").Value = -0.05309936136278375
$ws.Range("Z2").Value = 0.6208842897460018
$ws.Range("AA2").Value = 0.07791335411555654
$ws.Range("AB2").Value = 0.1351317619518854
$ws.Range("AC2").Value = -0.05721840783632882
$ws.Range("AG2").Value = -2.55
$ws.Range("AJ2").Value = -0.04623753399818675
$ws.Range("AK2").Value = -0.03583977512297962
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AO2").ClearContents()
$ws.Range("AP2").Value = -0.3935185185185185
$ws.Range("AQ2").ClearContents()

# --- Row 3 ---
$ws.Range("D3").Value = 0.117
$ws.Range("E3").Value = 0.0236
$ws.Range("G3").Value = 0.06515151515151515
$ws.Range("H3").Value = 0.06515151515151515
$ws.Range("I3").Value = 0.1608585858585859
$ws.Range("J3").Value = 0.1254877203406615
$ws.Range("K3").Value = 5.57
$ws.Range("L3").Value = 0.1406565656565656
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 2.55
$ws.Range("V3").Value = 0.04419410745233968
$ws.Range("W3").Value = 0.08203240058910162
$ws.Range("X3").Value = 0.1351317619518854
$ws.Range("Y3").Value = -0.05309936136278375
$ws.Range("Z3").Value = 0.6208842897460018
$ws.Range("AA3").Value = 0.07791335411555654
$ws.Range("AB3").Value = 0.1351317619518854
$ws.Range("AC3").Value = -0.05721840783632882
$ws.Range("AG3").Value = -2.55
$ws.Range("AJ3").Value = -0.04623753399818675
$ws.Range("AK3").Value = -0.03583977512297962
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AO3").ClearContents()
$ws.Range("AP3").Value = -0.3935185185185185
$ws.Range("AQ3").ClearContents()
